$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 4999
$ws.Range("I21").Value = 4999
$ws.Range("K21").Value = 4999
$ws.Range("M21").Value = -4531
$ws.Range("H23").Value = 4999
$ws.Range("I23").Value = 4999
$ws.Range("K23").Value = 4999
$ws.Range("M23").Value = -4765
$ws.Range("H29").Value = 71.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 675.4545000000001
$ws.Range("I38").Value = 675.4545000000001
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2026.3635
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1654.3635
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 5483.5557
$ws.Range("I58").Value = 1155.25
$ws.Range("J58").Value = 8946.200000000001
$ws.Range("K58").Value = 3465.75
$ws.Range("L58").Value = 26838.6
$ws.Range("M58").Value = -3315.75
$ws.Range("N58").Value = -27138.6
$ws.Range("H87").Value = 23250
$ws.Range("J87").Value = 21000
$ws.Range("L87").Value = 21000
$ws.Range("N87").Value = -23496
$ws.Range("H90").Value = 23250
$ws.Range("J90").Value = 21000
$ws.Range("L90").Value = 63000
$ws.Range("N90").Value = -75480
$ws.Range("H92").Value = 2042.9
$ws.Range("I92").Value = 388.06668
$ws.Range("K92").Value = 388.06668
$ws.Range("M92").Value = 859.93332
$ws.Range("H132").Value = 264440.9
$ws.Range("I132").Value = 1146.875
$ws.Range("J132").Value = 1668675.9
$ws.Range("K132").Value = 3440.625
$ws.Range("L132").Value = 5006027.699999999
$ws.Range("M132").Value = -910.625
$ws.Range("N132").Value = -5011087.699999999
$ws.Range("H135").Value = 10521.235
$ws.Range("I135").Value = 2135.0667
$ws.Range("K135").Value = 19215.6003
$ws.Range("M135").Value = -16680.6003
$ws.Range("H138").Value = 3758.4333
$ws.Range("I138").Value = 3722.3057
$ws.Range("J138").Value = 3782.5186
$ws.Range("K138").Value = 11166.9171
$ws.Range("L138").Value = 11347.5558
$ws.Range("M138").Value = -6026.917099999999
$ws.Range("N138").Value = -21627.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1475.7
$ws.Range("I2").Value = 1473.5883
$ws.Range("K2").Value = 1473.5883
$ws.Range("M2").Value = -1360.5883
$ws.Range("H31").Value = 1811.125
$ws.Range("I31").Value = 1811.125
$ws.Range("K31").Value = 1811.125
$ws.Range("M31").Value = -1517.125
$ws.Range("H116").Value = 1475.7
$ws.Range("I116").Value = 1473.5883
$ws.Range("K116").Value = 1473.5883
$ws.Range("M116").Value = 820.4117000000001
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H130").Value = 35833
$ws.Range("J130").Value = 35833
$ws.Range("L130").Value = 35833
$ws.Range("N130").Value = -45873

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1475.7
$ws.Range("I3").Value = 1473.5883
$ws.Range("K3").Value = 1473.5883
$ws.Range("M3").Value = -1359.5883
$ws.Range("H94").Value = 2280.484
$ws.Range("I94").Value = 775.9583
$ws.Range("J94").Value = 7438.857
$ws.Range("K94").Value = 775.9583
$ws.Range("L94").Value = 7438.857
$ws.Range("M94").Value = -324.9583
$ws.Range("N94").Value = -8340.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2311.7073
$ws.Range("I58").Value = 2141.4595
$ws.Range("J58").Value = 3886.5
$ws.Range("K58").Value = 2141.4595
$ws.Range("L58").Value = 3886.5
$ws.Range("M58").Value = -1938.4595
$ws.Range("N58").Value = -4292.5
$ws.Range("H105").Value = 4312.4707
$ws.Range("I105").Value = 4009.4167
$ws.Range("J105").Value = 5039.8
$ws.Range("K105").Value = 4009.4167
$ws.Range("L105").Value = 5039.8
$ws.Range("M105").Value = -2262.4167
$ws.Range("N105").Value = -8533.799999999999
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 29998.908
$ws.Range("I121").Value = 29998
$ws.Range("K121").Value = 29998
$ws.Range("M121").Value = -28688
$ws.Range("H136").Value = 2311.7073
$ws.Range("I136").Value = 2141.4595
$ws.Range("J136").Value = 3886.5
$ws.Range("K136").Value = 6424.3785
$ws.Range("L136").Value = 11659.5
$ws.Range("M136").Value = -3874.3785
$ws.Range("N136").Value = -16759.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6826.625
$ws.Range("I70").Value = 4942.5
$ws.Range("J70").Value = 9464.4
$ws.Range("K70").Value = 4942.5
$ws.Range("L70").Value = 9464.4
$ws.Range("M70").Value = -4672.5
$ws.Range("N70").Value = -10004.4
$ws.Range("H73").Value = 6826.625
$ws.Range("I73").Value = 4942.5
$ws.Range("J73").Value = 9464.4
$ws.Range("K73").Value = 4942.5
$ws.Range("L73").Value = 9464.4
$ws.Range("M73").Value = -4006.5
$ws.Range("N73").Value = -11336.4
$ws.Range("H132").Value = 2037.7273
$ws.Range("I132").Value = 2041.6
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 6124.799999999999
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -3594.799999999999
$ws.Range("N132").Value = -11057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1633.6471
$ws.Range("I22").Value = 1034.4445
$ws.Range("K22").Value = 1034.4445
$ws.Range("M22").Value = -739.4445000000001
$ws.Range("H27").Value = 1633.6471
$ws.Range("I27").Value = 1034.4445
$ws.Range("K27").Value = 1034.4445
$ws.Range("M27").Value = -927.4445000000001
$ws.Range("H46").Value = 2536.95
$ws.Range("I46").Value = 1793.9
$ws.Range("J46").Value = 3280
$ws.Range("K46").Value = 1793.9
$ws.Range("L46").Value = 3280
$ws.Range("M46").Value = -1605.9
$ws.Range("N46").Value = -3656
$ws.Range("H122").Value = 5149.143
$ws.Range("I122").Value = 4967.6523
$ws.Range("J122").Value = 5984
$ws.Range("K122").Value = 14902.9569
$ws.Range("L122").Value = 17952
$ws.Range("M122").Value = -12452.9569
$ws.Range("N122").Value = -22852
$ws.Range("H124:L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("N128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("H139:M139").ClearContents()
$ws.Range("H140:M140").ClearContents()
$ws.Range("H141:L141").ClearContents()
